# Add 2022-Q3 data:
#  1. Insert a brand-new worksheet named "2022-Q3" right before the existing
#     "2022-Q2" sheet (i.e. as the 2nd sheet, right after "总计"), and fill it
#     with the quarterly fund-holding detail rows.
#  2. Update the "总计" (summary) sheet: shift the existing quarter rows down
#     by one and put the new 2022-Q3 summary figures into row 2, with the
#     oldest quarter (2020-Q4) now landing on a new row 9.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3rows = @(
    @(0, "011834", "大成投资严选六月持有混合A", "3.10", "66.75", "6.17", "0.1913", 3),
    @(1, "013463", "大成致远优势一年持有期混合A", "3.65", "60.88", "3.57", "0.1303", 7),
    @(2, "011835", "大成投资严选六月持有混合C", "0.22", "66.75", "6.17", "0.0136", 3),
    @(3, "162416", "华宝港股通恒生香港35指数（LOF）", "0.20", "93.77", "4.14", "0.0083", 8),
    @(4, "013464", "大成致远优势一年持有期混合C", "0.17", "60.88", "3.57", "0.0061", 7)
)

$r = 2
foreach ($row in $q3rows) {
    $q3.Range("A$r").Value = $row[0]
    $q3.Range("B$r").Value = $row[1]
    $q3.Range("C$r").Value = $row[2]

    $q3.Range("D$r").NumberFormat = "@"
    $q3.Range("D$r").Value = $row[3]
    $q3.Range("E$r").NumberFormat = "@"
    $q3.Range("E$r").Value = $row[4]
    $q3.Range("F$r").NumberFormat = "@"
    $q3.Range("F$r").Value = $row[5]
    $q3.Range("G$r").NumberFormat = "@"
    $q3.Range("G$r").Value = $row[6]

    $q3.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update the "总计" summary sheet: shift rows 2-8 down to rows 3-9 (in
#    terms of content) and place the new 2022-Q3 totals on row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$totalRows = @(
    @("2022-Q3", 5, 0.35),
    @("2022-Q2", 6, 0.41),
    @("2022-Q1", 7, 0.46),
    @("2021-Q4", 7, 2.71),
    @("2021-Q3", 1, 0.01),
    @("2021-Q2", 3, 0.04),
    @("2021-Q1", 2, 0.02),
    @("2020-Q4", 1, 0.01)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

# Row 9 is brand new - give it the same index style/value as the rest of
# column A (sequential 0..7). Copy A8's formatting down first so the new
# cell picks up the same style as its neighbours, then set its value.
$total.Range("A8").Copy($total.Range("A9"))
$total.Range("A9").Value = 7
